$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Column D (Price) updates
    Set-TextValue $ws "D2" '56.420.80'
    Set-TextValue $ws "D3" '2.941.47'
    Set-TextValue $ws "D5" '493.02'
    Set-TextValue $ws "D6" '133.32'
    Set-TextValue $ws "D9" '7.10'
    Set-TextValue $ws "D11" '0.350'
    Set-TextValue $ws "D12" '3.446.90'
    Set-TextValue $ws "D14" '25.85'
    Set-TextValue $ws "D15" '0.0000156'
    Set-TextValue $ws "D16" '56.498.05'
    Set-TextValue $ws "D17" '5.95'
    Set-TextValue $ws "D18" '2.942.79'
    Set-TextValue $ws "D19" '12.42'
    Set-TextValue $ws "D21" '315.20'
    Set-TextValue $ws "D23" '5.75'
    Set-TextValue $ws "D24" '0.482'
    Set-TextValue $ws "D25" '62.40'
    Set-TextValue $ws "D26" '0.997'
    Set-TextValue $ws "D29" '6.41'
    Set-TextValue $ws "D30" '6.98'
    Set-TextValue $ws "D31" '1.73'
    Set-TextValue $ws "D32" '19.86'
    Set-TextValue $ws "D33" '1.12'
    Set-TextValue $ws "D34" '150.97'
    Set-TextValue $ws "D35" '4.42'
    Set-TextValue $ws "D36" '5.64'
    Set-TextValue $ws "D38" '23.67'
    Set-TextValue $ws "D39" '0.0649'
    Set-TextValue $ws "D41" '2.970.40'
    Set-TextValue $ws "D43" '3.66'
    Set-TextValue $ws "D44" '0.634'
    Set-TextValue $ws "D45" '2.131.44'
    Set-TextValue $ws "D46" '1.33'
    Set-TextValue $ws "D47" '5.82'
    Set-TextValue $ws "D48" '0.904'
    Set-TextValue $ws "D49" '0.0228'
    Set-TextValue $ws "D50" '18.79'
    Set-TextValue $ws "D51" '0.0845'

# Column E (Volume(1h)) updates
    $ws.Range("E2").Value = '  -3.40%  '
    $ws.Range("E3").Value = '  -4.32%  '
    $ws.Range("E4").Value = '  -0.02%  '
    $ws.Range("E5").Value = '  -6.97%  '
    $ws.Range("E6").Value = '  -7.24%  '
    $ws.Range("E7").Value = '  +0.16%  '
    $ws.Range("E8").Value = '  -5.94%  '
    $ws.Range("E9").Value = '  -6.94%  '
    $ws.Range("E10").Value = '  -7.76%  '
    $ws.Range("E11").Value = '  -6.07%  '
    $ws.Range("E12").Value = '  -4.19%  '
    $ws.Range("E13").Value = '  -3.78%  '
    $ws.Range("E14").Value = '  -6.00%  '
    $ws.Range("E15").Value = '  -10.73%  '
    $ws.Range("E16").Value = '  -3.16%  '
    $ws.Range("E17").Value = '  -4.68%  '
    $ws.Range("E18").Value = '  -4.39%  '
    $ws.Range("E19").Value = '  -6.04%  '
    $ws.Range("E20").Value = '  -6.22%  '
    $ws.Range("E21").Value = '  -8.29%  '
    $ws.Range("E22").Value = '  -0.11%  '
    $ws.Range("E23").Value = '  -0.10%  '
    $ws.Range("E24").Value = '  -5.14%  '
    $ws.Range("E25").Value = '  -4.93%  '
    $ws.Range("E26").Value = '  -0.07%  '
    $ws.Range("E27").Value = '  -4.99%  '
    $ws.Range("E28").Value = '  -12.81%  '
    $ws.Range("E29").Value = '  -9.11%  '
    $ws.Range("E30").Value = '  -7.36%  '
    $ws.Range("E31").Value = '  -7.33%  '
    $ws.Range("E32").Value = '  -6.64%  '
    $ws.Range("E33").Value = '  -9.61%  '
    $ws.Range("E34").Value = '  -4.73%  '
    $ws.Range("E35").Value = '  -8.78%  '
    $ws.Range("E36").Value = '  -6.11%  '
    $ws.Range("E37").Value = '  -10.30%  '
    $ws.Range("E38").Value = '  -10.14%  '
    $ws.Range("E39").Value = '  -7.62%  '
    $ws.Range("E40").Value = '  -1.36%  '
    $ws.Range("E41").Value = '  -4.43%  '
    $ws.Range("E42").Value = '  -0.05%  '
    $ws.Range("E43").Value = '  -7.89%  '
    $ws.Range("E44").Value = '  -5.38%  '
    $ws.Range("E45").Value = '  -8.97%  '
    $ws.Range("E46").Value = '  -10.25%  '
    $ws.Range("E47").Value = '  -4.33%  '
    $ws.Range("E48").Value = '  -13.61%  '
    $ws.Range("E49").Value = '  -6.57%  '
    $ws.Range("E50").Value = '  -7.22%  '
    $ws.Range("E51").Value = '  -6.96%  '
